$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing data shifts down from row 1 to row 2, etc.
$ws.Rows.Item(1).Insert()

# Set the new title cell and merge A1:H1
$ws.Range("A1").Value = "1st Semester"
$ws.Range("A1:H1").Merge()
$ws.Range("A1:H1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A1:H1").Borders.Item(9).LineStyle = 1    # xlEdgeBottom, xlContinuous
$ws.Range("A1:H1").Borders.Item(9).Weight = 2       # xlThin
